$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1833.0834
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 1599.4
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 1599.4
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -1737.4
$ws.Range("H112").Value = 1971.8422
$ws.Range("J112").Value = 1971.8422
$ws.Range("L112").Value = 5915.5266
$ws.Range("N112").Value = -8131.5266
$ws.Range("H116").Value = 16893.875
$ws.Range("J116").Value = 25001
$ws.Range("L116").Value = 25001
$ws.Range("N116").Value = -31885
$ws.Range("H125").Value = 3079.6
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
$ws.Range("H135").Value = 1546.0667
$ws.Range("I135").Value = 1370.8572
$ws.Range("J135").Value = 3999
$ws.Range("K135").Value = 12337.7148
$ws.Range("L135").Value = 35991
$ws.Range("M135").Value = -9802.7148
$ws.Range("N135").Value = -41061
$ws.Range("H137").Value = 5518.273
$ws.Range("I137").Value = 4753.615
$ws.Range("J137").Value = 6622.778
$ws.Range("K137").Value = 14260.845
$ws.Range("L137").Value = 19868.334
$ws.Range("M137").Value = -11710.845
$ws.Range("N137").Value = -24968.334
$ws.Range("H138").Value = 5147.6
$ws.Range("I138").Value = 3745.8125
$ws.Range("J138").Value = 5414.607
$ws.Range("K138").Value = 11237.4375
$ws.Range("L138").Value = 16243.821
$ws.Range("M138").Value = -6097.4375
$ws.Range("N138").Value = -26523.821

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12634.88
$ws.Range("I32").Value = 6197.0283
$ws.Range("J32").Value = 28396.518
$ws.Range("K32").Value = 6197.0283
$ws.Range("L32").Value = 28396.518
$ws.Range("M32").Value = -5910.0283
$ws.Range("N32").Value = -28970.518
$ws.Range("H46").Value = 2019.75
$ws.Range("I46").Value = 1656.3334
$ws.Range("K46").Value = 1656.3334
$ws.Range("M46").Value = -1337.3334
$ws.Range("H61").Value = 4363.064
$ws.Range("I61").Value = 3933.6191
$ws.Range("J61").Value = 7970.4
$ws.Range("K61").Value = 3933.6191
$ws.Range("L61").Value = 7970.4
$ws.Range("M61").Value = -3721.6191
$ws.Range("N61").Value = -8394.4
$ws.Range("H63").Value = 6429.0835
$ws.Range("I63").Value = 3999.8333
$ws.Range("J63").Value = 8858.333000000001
$ws.Range("K63").Value = 3999.8333
$ws.Range("L63").Value = 8858.333000000001
$ws.Range("M63").Value = -3313.8333
$ws.Range("N63").Value = -10230.333
$ws.Range("H66").Value = 6429.0835
$ws.Range("I66").Value = 3999.8333
$ws.Range("J66").Value = 8858.333000000001
$ws.Range("K66").Value = 19999.1665
$ws.Range("L66").Value = 44291.665
$ws.Range("M66").Value = -16567.1665
$ws.Range("N66").Value = -51155.665
$ws.Range("H94").Value = 29330
$ws.Range("J94").Value = 29330
$ws.Range("L94").Value = 29330
$ws.Range("N94").Value = -31132
$ws.Range("H132").Value = 4636.4062
$ws.Range("I132").Value = 3526.9363
$ws.Range("K132").Value = 10580.8089
$ws.Range("M132").Value = -8050.8089
$ws.Range("H136").Value = 4363.064
$ws.Range("I136").Value = 3933.6191
$ws.Range("J136").Value = 7970.4
$ws.Range("K136").Value = 11800.8573
$ws.Range("L136").Value = 23911.2
$ws.Range("M136").Value = -9250.8573
$ws.Range("N136").Value = -29011.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27981.977
$ws.Range("I134").Value = 5102.61
$ws.Range("K134").Value = 15307.83
$ws.Range("M134").Value = -12772.83

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1977.625
$ws.Range("I31").Value = 1570.1666
$ws.Range("K31").Value = 1570.1666
$ws.Range("M31").Value = -1275.1666
$ws.Range("H34").Value = 1977.625
$ws.Range("I34").Value = 1570.1666
$ws.Range("K34").Value = 1570.1666
$ws.Range("M34").Value = -1368.1666
$ws.Range("H58").Value = 1769.12
$ws.Range("I58").Value = 1473.9231
$ws.Range("J58").Value = 2088.9167
$ws.Range("K58").Value = 1473.9231
$ws.Range("L58").Value = 2088.9167
$ws.Range("M58").Value = -1270.9231
$ws.Range("N58").Value = -2494.9167
$ws.Range("H99").Value = 5878.25
$ws.Range("I99").Value = 6753
$ws.Range("J99").Value = 5003.5
$ws.Range("K99").Value = 6753
$ws.Range("L99").Value = 5003.5
$ws.Range("M99").Value = -5255
$ws.Range("N99").Value = -7999.5
$ws.Range("H126").Value = 5878.25
$ws.Range("I126").Value = 6753
$ws.Range("J126").Value = 5003.5
$ws.Range("K126").Value = 20259
$ws.Range("L126").Value = 15010.5
$ws.Range("M126").Value = -17789
$ws.Range("N126").Value = -19950.5
$ws.Range("H130").Value = 86998
$ws.Range("J130").Value = 86998
$ws.Range("L130").Value = 86998
$ws.Range("N130").Value = -97038
$ws.Range("H132").Value = 2621.0789
$ws.Range("I132").Value = 2189.84
$ws.Range("J132").Value = 3450.3845
$ws.Range("K132").Value = 6569.52
$ws.Range("L132").Value = 10351.1535
$ws.Range("M132").Value = -4039.52
$ws.Range("N132").Value = -15411.1535
$ws.Range("H136").Value = 1769.12
$ws.Range("I136").Value = 1473.9231
$ws.Range("J136").Value = 2088.9167
$ws.Range("K136").Value = 4421.7693
$ws.Range("L136").Value = 6266.750100000001
$ws.Range("M136").Value = -1871.7693
$ws.Range("N136").Value = -11366.7501
$ws.Range("H141").Value = 509763.4
$ws.Range("J141").Value = 540762.25
$ws.Range("L141").Value = 540762.25
$ws.Range("N141").Value = -551122.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 211.95653
$ws.Range("I2").Value = 88.44444
$ws.Range("J2").Value = 291.35715
$ws.Range("K2").Value = 530.66664
$ws.Range("L2").Value = 1748.1429
$ws.Range("M2").Value = -417.66664
$ws.Range("N2").Value = -1974.1429
$ws.Range("H40").Value = 175.7
$ws.Range("I40").Value = 101.875
$ws.Range("K40").Value = 407.5
$ws.Range("M40").Value = -338.5
$ws.Range("H76").Value = 5007.5
$ws.Range("J76").Value = 5015
$ws.Range("L76").Value = 15045
$ws.Range("N76").Value = -15811
$ws.Range("H79").Value = 5007.5
$ws.Range("J79").Value = 5015
$ws.Range("L79").Value = 15045
$ws.Range("N79").Value = -17697
$ws.Range("H108").Value = 12121
$ws.Range("I108").Value = 12121
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 36363
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -33483
$ws.Range("H120").Value = 268332.75
$ws.Range("I120").Value = 504999
$ws.Range("J120").Value = 31666.5
$ws.Range("K120").Value = 1514997
$ws.Range("L120").Value = 94999.5
$ws.Range("M120").Value = -1510159
$ws.Range("N120").Value = -104675.5
$ws.Range("H122").Value = 1377
$ws.Range("J122").Value = 1498.091
$ws.Range("L122").Value = 13482.819
$ws.Range("N122").Value = -18382.819
$ws.Range("H129").Value = 47763104
$ws.Range("J129").Value = 502250
$ws.Range("L129").Value = 1506750
$ws.Range("N129").Value = -1516750
$ws.Range("H133").Value = 19912.834
$ws.Range("I133").Value = 19895.4
$ws.Range("J133").Value = 20000
$ws.Range("K133").Value = 59686.2
$ws.Range("L133").Value = 60000
$ws.Range("M133").Value = -54626.2
$ws.Range("N133").Value = -70120
$ws.Range("H138").Value = 4990.75
$ws.Range("J138").Value = 14989.5
$ws.Range("L138").Value = 44968.5
$ws.Range("N138").Value = -55248.5
$ws.Range("N108").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I70").Value = 7566.7144
$ws.Range("K70").Value = 7566.7144
$ws.Range("M70").Value = -7296.7144
$ws.Range("I73").Value = 7566.7144
$ws.Range("K73").Value = 7566.7144
$ws.Range("M73").Value = -6630.7144
$ws.Range("H126").Value = 4857
$ws.Range("J126").Value = 4857
$ws.Range("L126").Value = 14571
$ws.Range("N126").Value = -19511
$ws.Range("H132").Value = 27456.113
$ws.Range("I132").Value = 4867.2
$ws.Range("J132").Value = 75860.92999999999
$ws.Range("K132").Value = 14601.6
$ws.Range("L132").Value = 227582.79
$ws.Range("M132").Value = -12071.6
$ws.Range("N132").Value = -232642.79

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6019.125
$ws.Range("I7").Value = 4000
$ws.Range("K7").Value = 4000
$ws.Range("M7").Value = -3888
$ws.Range("H126").Value = 6019.125
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530
$ws.Range("H132").Value = 5064.577
$ws.Range("I132").Value = 4951.7827
$ws.Range("K132").Value = 14855.3481
$ws.Range("M132").Value = -12325.3481

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1857.5714
$ws.Range("J107").Value = 1000.5
$ws.Range("L107").Value = 3001.5
$ws.Range("N107").Value = -6841.5
$ws.Range("H132").Value = 31925.883
$ws.Range("I132").Value = 2276.125
$ws.Range("J132").Value = 103085.3
$ws.Range("K132").Value = 6828.375
$ws.Range("L132").Value = 309255.9
$ws.Range("M132").Value = -4298.375
$ws.Range("N132").Value = -314315.9
$ws.Range("H136").Value = 261994.88
$ws.Range("I136").Value = 250993.83
$ws.Range("K136").Value = 752981.49
$ws.Range("M136").Value = -750431.49

